# Update column G ("K") values per regenerated save_data
# (source data changed to compute K from "K" field instead of "Strike#";
#  std/mean were recalculated upstream and these are the refreshed values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 2
    9  = 0
    10 = 1
    11 = 0
    12 = 3
    14 = 2
    15 = 0
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 2
    22 = 0
    23 = 1
    24 = 2
    25 = 1
    26 = 3
    28 = 1
    29 = 2
    30 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
